$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (columns K through T)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.367603
$ws.Range("N2").Value = 1.102809
$ws.Range("O2").Value = 0.5971364972068339
$ws.Range("P2").Value = 0.5971364972068339
$ws.Range("Q2").Value = 0.1668459341593333
$ws.Range("R2").Value = 1.501613407434
$ws.Range("S2").Value = 0.5971364972068339
$ws.Range("T2").Value = 0.5971364972068339

# Update row 3 values (columns O, P, S, T)
$ws.Range("O3").Value = 0.4028635027931661
$ws.Range("P3").Value = 0.402863502793166
$ws.Range("S3").Value = 0.4028635027931661
$ws.Range("T3").Value = 0.402863502793166

# Delete row 4 entirely (it referenced the "MuSCs" target cluster)
$ws.Range("A4:T4").EntireRow.Delete()
